$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.407.18'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '2.013.64'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.41'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.77%  '

$ws.Range("E9").Value = '  +1.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0771'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.25%  '

$ws.Range("E11").Value = '  -2.08%  '

$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.65%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.308.74'
$ws.Range("E13").Value = '  +0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.807'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.40%  '

$ws.Range("E16").Value = '  -3.03%  '

$ws.Range("D17").Value = '2.017.38'
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").Value = '37.259.39'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '

$ws.Range("D20").Value = '0.0₃0840'
$ws.Range("E20").Value = '  -2.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.54%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.49%  '

$ws.Range("E29").Value = '  -9.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.59%  '

$ws.Range("E31").Value = '  -0.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0650'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.82%  '

$ws.Range("E36").Value = '  +0.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.98%  '

$ws.Range("E40").Value = '  +3.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.33%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0213'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.38%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0931'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.10%  '

$ws.Range("D44").Value = '1.392.41'
$ws.Range("E44").Value = '  +1.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.50%  '

$ws.Range("E47").Value = '  -1.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.03%  '

$ws.Range("E49").Value = '  +1.89%  '

$ws.Range("D50").Value = '2.200.39'

$ws.Range("E51").Value = '  -3.69%  '
